$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (and one coin swap in row 51).
# NumberFormat is forced to text ("@") per cell before assigning the new
# value so numeric-looking strings (e.g. "216.58") are preserved as text,
# matching how the source data is stored, instead of Excel auto-converting
# them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.001.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.186.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.58'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '633.72'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.396'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.721'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +6.17%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.184.23'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.569'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.183'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.808.68'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.35'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.775.81'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '32.63'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.192.24'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000217'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +47.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.48'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '435.49'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.45%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.99'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.29'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.66'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '80.95'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +10.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.347.85'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.159'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.04'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +29.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.37'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '513.90'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.99'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.31%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.29'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.34'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.40'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.127'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.47%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.373'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '146.30'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '44.01'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '169.34'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.124'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.740'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.87'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.04%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.21'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.12%  '
